$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.981.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "'1.634.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'212.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.64%  "

$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").Value = "'23.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.35%  "

$ws.Range("E9").Value = "  -2.24%  "

$ws.Range("E11").Value = "  +0.84%  "

$ws.Range("D12").Value = "'1.865.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("D13").Value = "'1.634.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("E14").Value = "  -0.30%  "

$ws.Range("D15").Value = "'0.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.97%  "

$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("D17").Value = "'27.969.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").Value = "'232.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.73%  "

$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("E20").Value = "  -0.80%  "

$ws.Range("E21").Value = "  -0.29%  "

$ws.Range("E22").Value = "  -4.37%  "

$ws.Range("E23").Value = "  -0.84%  "

$ws.Range("E24").Value = "  -3.37%  "

$ws.Range("D25").Value = "'154.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "

$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("E28").Value = "  -0.42%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("E31").Value = "  -0.71%  "

$ws.Range("E32").Value = "  +2.15%  "

$ws.Range("D33").Value = "'3.10"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.14%  "

$ws.Range("D34").Value = "'1.410.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.10%  "

$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.22%  "

$ws.Range("E37").Value = "  +0.57%  "

$ws.Range("E38").Value = "  +1.93%  "

$ws.Range("E39").Value = "  +0.18%  "

$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("E41").Value = "  -1.10%  "

$ws.Range("D43").Value = "'67.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.96%  "

$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("E45").Value = "  +0.82%  "

$ws.Range("D46").Value = "'2.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").Value = "'1.775.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").Value = "'88.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.19%  "

$ws.Range("D49").Value = "'0.1000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "

$ws.Range("E50").Value = "  -0.24%  "

$ws.Range("D51").Value = "'7.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.91%  "
